# Swap the B:AD content between pairs of adjacent data rows.
# Column A (the sequential row index) is left untouched on each row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(215, 216),
    @(229, 231),
    @(232, 233),
    @(245, 246),
    @(248, 249),
    @(251, 252),
    @(271, 272),
    @(310, 311),
    @(316, 317)
)

foreach ($p in $pairs) {
    $r1 = $p[0]
    $r2 = $p[1]
    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}
